$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-21 Saturday" "2024-12-22 Sunday"

Replace-Text "41×40=" "41×32="
Replace-Text "13×17=" "32×81="
Replace-Text "77×72=" "36×33="
Replace-Text "19×34=" "91×71="
Replace-Text "62×86=" "44×95="

Replace-Text "77×79=" "66×63="
Replace-Text "63×91=" "56×26="
Replace-Text "85×62=" "91×97="
Replace-Text "83×75=" "26×34="
Replace-Text "87×59=" "53×37="

Replace-Text "81×15=" "61×13="
Replace-Text "35×72=" "79×44="
Replace-Text "59×55=" "16×85="
Replace-Text "62×46=" "99×95="
Replace-Text "50×34=" "27×81="

Replace-Text "31×75=" "72×53="
Replace-Text "41×25=" "13×89="
Replace-Text "41×71=" "83×89="
Replace-Text "41×11=" "84×91="
Replace-Text "79×52=" "78×34="

Replace-Text "90×36=" "11×49="
Replace-Text "20×33=" "81×31="
Replace-Text "90×41=" "53×60="
Replace-Text "45×53=" "17×85="
Replace-Text "88×63=" "30×65="
